$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.885.90"
$ws.Cells.Item(2, 5).Value = "  -0.97%  "

$ws.Cells.Item(3, 4).Value = "1.638.31"
$ws.Cells.Item(3, 5).Value = "  -1.28%  "

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.ClearFormats()
$ws.Cells.Item(4, 5).Value = "  -0.09%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "216.22"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -1.02%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5038"
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -2.12%  "

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

$ws.Cells.Item(8, 2).Value = "Cardano"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2569"
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -0.61%  "

$ws.Cells.Item(9, 2).Value = "Dogecoin"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06431"
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +0.05%  "

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.58"
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = "  -1.73%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07702"
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -1.08%  "

$ws.Cells.Item(12, 4).Value = "1.646.59"
$ws.Cells.Item(12, 5).Value = "  -1.05%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.245"
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -1.18%  "

$ws.Cells.Item(14, 4).Value = "1.867.88"
$ws.Cells.Item(14, 5).Value = "  -1.08%  "

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5461"
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = "  -1.68%  "

$ws.Cells.Item(16, 4).Value = "0.0₅7946"
$ws.Cells.Item(16, 5).Value = "  -1.23%  "

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "63.47"
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  -1.24%  "

$ws.Cells.Item(18, 4).Value = "25.933.55"
$ws.Cells.Item(18, 5).Value = "  -1.03%  "

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -0.10%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "203.34"
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -4.03%  "

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.310"
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  -2.33%  "

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.974"
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -0.66%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.993"
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -0.01%  "

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +0.00%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.954"
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +11.04%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "141.73"
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -1.59%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1150"
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  -1.11%  "

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.72"
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -0.57%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.720"
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = "  -3.72%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05038"
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -4.55%  "

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.242"
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  -1.04%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.259"
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -3.08%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.189"
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -0.85%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.535"
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -2.32%  "

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.357"
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = "  -0.66%  "

$ws.Cells.Item(36, 4).Value = "1.167.60"
$ws.Cells.Item(36, 5).Value = "  -0.11%  "

$ws.Cells.Item(37, 2).Value = "MXToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.633"
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -4.68%  "

$ws.Cells.Item(38, 2).Value = "ARBITRUM"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.8929"
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -3.75%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5559"
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -1.69%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.01563"
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = "  -1.74%  "

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.565"
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = "  -0.11%  "

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -0.05%  "

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.663"
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -0.60%  "

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.8075"
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = "  -4.02%  "

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "99.76"
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -0.54%  "

$ws.Cells.Item(46, 4).Value = "1.780.45"
$ws.Cells.Item(46, 5).Value = "  -0.98%  "

$ws.Cells.Item(47, 4).Value = "0.0₈109"
$ws.Cells.Item(47, 5).Value = "  -4.39%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.4522"
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  -0.35%  "

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.011"
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +0.49%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "55.04"
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = "  -1.51%  "

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05038"
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -0.34%  "
